# The "ja" sheet lists resource-bundle rows in column A using a running
# index formula (A17:A81 is one shared-formula group: each cell = previous+1,
# A82:A87 is a second shared-formula group).
#
# Row 56 (A56) was a manual/independent override "=A54+1" that evaluated to
# 41 (duplicating the index already used by A55). The edit turns it into a
# normal "+1 from the row above" formula ("=A55+1", matching the pattern of
# the surrounding shared formula group), which bumps its own value to 42 and
# cascades a +1 shift through every following row's computed index
# (43, 44, 45, ... 73) all the way down to row 87.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A56").Formula = "=A55+1"

# Update the view: scroll position moved down and the active selection
# moved to A90 (just below the data, which ends at row 88).
$excel.ActiveWindow.ScrollRow = 54
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A90").Select()
